$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "04:03:59.302590400"
$ws.Range("B3").Value = "saucedemo2"
$ws.Range("C3").Value = "04:04:36.417331300"
$ws.Range("B4").Value = "saucedemo1"
$ws.Range("C4").Value = "04:04:36.986213500"
